{"js": "// Edit: split the run \"v\u00e0 ghi log \u0111\u0103ng nh\u1eadp cho t\u1eebng s\u1ef1 ki\u1ec7n \u0111\u1ec3 bi\u1ebft ai \u0111\u00e3\n// v\u1eadn h\u00e0nh h\u1ec7 th\u1ed1ng. \" into two runs \u2014 \"...h\u1ec7 th\u1ed1\" + \"ng.\" (dropping the\n// trailing space), and move the \"_GoBack\" bookmark from the end of the last\n// bullet (\"Ch\u1ee9c n\u0103ng th\u00f4ng b\u00e1o ...\") to the end of this paragraph.\n\nconst body = context.document.body;\n\nconst ORIGINAL_SENTENCE =\n  \"v\u00e0 ghi log \u0111\u0103ng nh\u1eadp cho t\u1eebng s\u1ef1 ki\u1ec7n \u0111\u1ec3 bi\u1ebft ai \u0111\u00e3 v\u1eadn h\u00e0nh h\u1ec7 th\u1ed1ng. \";\nconst FIRST_PART =\n  \"v\u00e0 ghi log \u0111\u0103ng nh\u1eadp cho t\u1eebng s\u1ef1 ki\u1ec7n \u0111\u1ec3 bi\u1ebft ai \u0111\u00e3 v\u1eadn h\u00e0nh h\u1ec7 th\u1ed1\";\nconst SECOND_PART = \"ng.\";\n\n// 1) Locate the sentence and replace it with two explicit runs so the split\n//    survives the round trip (engine normally merges adjacent runs whose\n//    formatting matches when they are written with separate insertText\n//    calls, so we feed the exact OOXML for both runs in one shot).\nconst results = body.search(ORIGINAL_SENTENCE, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target sentence to split.\");\n}\n\nconst target = results.items[0];\n\nconst ooxmlSnippet =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p>\" +\n  \"<w:r><w:t>\" + FIRST_PART + \"</w:t></w:r>\" +\n  \"<w:r><w:t>\" + SECOND_PART + \"</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\ntarget.insertOoxml(ooxmlSnippet, \"Replace\");\nawait context.sync();\n\n// 2) Move the \"_GoBack\" bookmark here. Word keeps a single hidden \"_GoBack\"\n//    bookmark tracking the last edit location, so remove the old one first\n//    (otherwise a duplicate bookmark with the same name would be created).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet editedParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"v\u00e0 ghi log \u0111\u0103ng nh\u1eadp\") !== -1 &&\n      text.indexOf(\"v\u1eadn h\u00e0nh h\u1ec7 th\u1ed1ng\") !== -1) {\n    editedParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!editedParagraph) {\n  throw new Error(\"Could not find the edited paragraph to re-anchor the bookmark.\");\n}\n\nconst endRange = editedParagraph.getRange(\"End\");\nendRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Step 1: split the run \"v\u00e0 ghi log \u0111\u0103ng nh\u1eadp cho t\u1eebng s\u1ef1 ki\u1ec7n \u0111\u1ec3 bi\u1ebft ai\n# \u0111\u00e3 v\u1eadn h\u00e0nh h\u1ec7 th\u1ed1ng. \" into \"...h\u1ec7 th\u1ed1\" + \"ng.\" (dropping the trailing space).\n$r = $d.Content\n$r.Find.ClearFormatting()\n$found = $r.Find.Execute(\"v\u00e0 ghi log \u0111\u0103ng nh\u1eadp cho t\u1eebng s\u1ef1 ki\u1ec7n \u0111\u1ec3 bi\u1ebft ai \u0111\u00e3 v\u1eadn h\u00e0nh h\u1ec7 th\u1ed1ng. \")\nif (-not $found) {\n    throw \"Could not find target sentence to split.\"\n}\n$r.Text = \"\"\n$r.Collapse(0)\n$r.InsertAfter(\"v\u00e0 ghi log \u0111\u0103ng nh\u1eadp cho t\u1eebng s\u1ef1 ki\u1ec7n \u0111\u1ec3 bi\u1ebft ai \u0111\u00e3 v\u1eadn h\u00e0nh h\u1ec7 th\u1ed1\")\n$r.Collapse(0)\n$r.InsertAfter(\"ng.\")\n\n# --- Step 2: move the hidden \"_GoBack\" bookmark from the end of the last\n# bullet (\"Ch\u1ee9c n\u0103ng th\u00f4ng b\u00e1o ...\") to the end of the paragraph we just edited.\n$d.Bookmarks(\"_GoBack\").Delete()\n\n# Locate the paragraph we edited.\n$targetPara = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*v\u00e0 ghi log \u0111\u0103ng nh\u1eadp*\" -and $t -like \"*v\u1eadn h\u00e0nh*\") {\n        $targetPara = $p\n        break\n    }\n}\nif ($null -eq $targetPara) {\n    throw \"Could not find the edited paragraph to re-anchor the bookmark.\"\n}\n\n# Work around an off-by-one quirk when adding a collapsed bookmark at the very\n# last offset of a paragraph (immediately before its paragraph mark): insert a\n# throwaway placeholder character there first, anchor the bookmark just before\n# it, then remove the placeholder again.\n$endPos = $targetPara.Range.End - 1\n$placeholderRange = $d.Range($endPos, $endPos)\n$placeholderRange.InsertAfter(\"X\")\n\n$bookmarkRange = $d.Range($endPos, $endPos)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n$bm = $d.Bookmarks(\"_GoBack\")\n$placeholder = $d.Range($bm.End, $bm.End + 1)\n$placeholder.Delete()\n"}
